# edit.ps1 -- applies the "typo fixes / City Storage Systems" revision to
# mslaterResume-W.docx
#
# Summary of changes (from the commit's unified diff):
#   1. The two header hyperlink runs (email + Portfolio) pick up a
#      character style reference (rStyle = "Style9"). That style does not
#      exist yet in styles.xml, so it has to be minted first.
#   2. "for the for the Otter and CloudKitchens" -> "for the Otter and
#      CloudKitchens" (duplicated words typo).
#   3. "ExactTarget" -> "Exact Target" (missing space).
#   4. "Wordpress" -> "WordPress" -- but landed as three separate runs
#      with identical run formatting (Photoshop, Word | P | ress, HTML...),
#      which is what you get from retyping just the "p" in place.
#   5. The Caption paragraph style's display name changes from "Caption"
#      to "caption" (style id is untouched).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Mint the "Style9" character style and apply it to both hyperlink runs
#    (mslater220@gmail.com, Portfolio) in the header.
# ---------------------------------------------------------------------
$d.Styles.Add("Style9", 2) | Out-Null

for ($i = 1; $i -le $d.Hyperlinks.Count; $i++) {
    $h = $d.Hyperlinks.Item($i)
    $h.Range.Style = "Style9"
}

# ---------------------------------------------------------------------
# 2. Fix the duplicated "for the for the" typo in the Otter/CloudKitchens
#    KB-article bullet.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Developed KB article content for the for the Otter and CloudKitchens brands.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Developed KB article content for the Otter and CloudKitchens brands.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 3. "ExactTarget" -> "Exact Target"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "ExactTarget", $true, $false, $false, $false, $false, $true, 1, $false,
    "Exact Target", 2) | Out-Null

# ---------------------------------------------------------------------
# 4. "Wordpress" -> "WordPress", retyped as three runs that share identical
#    formatting: "...Photoshop, Word" | "P" | "ress, HTML, CSS...".
#    A straight Find/Replace collapses back into a single run because the
#    formatting is unchanged, so the capital P is re-typed in place and a
#    harmless Bold on/off toggle is used to force the run split that a real
#    retype would have produced, without altering the visible formatting.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Wordpress", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $rng.Start
    $pRange = $d.Range($start + 4, $start + 5)
    $pRange.Text = "P"
    $pRange.Bold = 1
    $pRange.Bold = 0
}

# ---------------------------------------------------------------------
# 5. Rename the Caption style's display name to lowercase "caption"
#    (styleId stays "Caption").
# ---------------------------------------------------------------------
$capStyle = $d.Styles.Item("Caption")
$capStyle.NameLocal = "caption"

Write-Output "edit complete"
